$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.521.26"
$ws.Range("E2").Value = "  +0.87%  "

# Row 3
$ws.Range("D3").Value = "1.728.98"
$ws.Range("E3").Value = "  +0.42%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.34"
$ws.Range("E5").Value = "  +2.09%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  +0.02%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4814"
$ws.Range("E7").Value = "  +1.64%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2676"
$ws.Range("E8").Value = "  +1.63%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06194"
$ws.Range("E9").Value = "  -0.05%  "

# Row 10
$ws.Range("D10").Value = "1.731.03"
$ws.Range("E10").Value = "  +0.63%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07195"
$ws.Range("E11").Value = "  +1.94%  "

# Row 12
$ws.Range("E12").Value = "  +0.56%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6110"
$ws.Range("E13").Value = "  +1.92%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.530"
$ws.Range("E14").Value = "  +2.04%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.22"
$ws.Range("E15").Value = "  +1.22%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9997"
$ws.Range("E16").Value = "  +0.01%  "

# Row 17
$ws.Range("D17").Value = "26.516.46"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9998"
$ws.Range("E18").Value = "  +0.01%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006955"
$ws.Range("E19").Value = "  +2.00%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.56"
$ws.Range("E20").Value = "  +0.16%  "

# Row 21
$ws.Range("D21").Value = "1.954.18"
$ws.Range("E21").Value = "  +0.90%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.527"
$ws.Range("E22").Value = "  -0.22%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.819"
$ws.Range("E23").Value = "  +1.03%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.259"
$ws.Range("E24").Value = "  +0.09%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.83"
$ws.Range("E25").Value = "  +1.24%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.35"
$ws.Range("E26").Value = "  +0.88%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.775"
$ws.Range("E27").Value = "  +0.35%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.407"
$ws.Range("E28").Value = "  +0.50%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "107.57"
$ws.Range("E29").Value = "  +0.78%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.979"
$ws.Range("E30").Value = "  +0.58%  "

# Row 31
$ws.Range("E31").Value = "  +3.00%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.697"
$ws.Range("E32").Value = "  +0.22%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04521"
$ws.Range("E33").Value = "  +0.20%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.617"
$ws.Range("E34").Value = "  +0.11%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.003"
$ws.Range("E35").Value = "  +2.07%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6262"
$ws.Range("E36").Value = "  +0.34%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.078"
$ws.Range("E37").Value = "  +7.64%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9111"
$ws.Range("E38").Value = "  -2.41%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.387"
$ws.Range("E39").Value = "  -2.75%  "

# Row 41
$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "103.42"
$ws.Range("E41").Value = "  -9.97%  "

# Row 42
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01504"
$ws.Range("E42").Value = "  +1.24%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.539"
$ws.Range("E43").Value = "  -2.09%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3882"
$ws.Range("E44").Value = "  +1.06%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.989"
$ws.Range("E45").Value = "  +9.79%  "

# Row 46
$ws.Range("E46").Value = "  -0.38%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05365"
$ws.Range("E47").Value = "  +1.83%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.54"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.790"
$ws.Range("E49").Value = "  -1.36%  "

# Row 50
$ws.Range("E50").Value = "  +2.83%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3405"
$ws.Range("E51").Value = "  +0.43%  "
